# Update cryptocurrency price/volume data per Jun 3 2024 GitHub Actions scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '68.651.96'
$ws.Range('E2').Value = '  +1.25%  '
$ws.Range('D3').Value = '3.814.94'
$ws.Range('E3').Value = '  +0.21%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '613.16'
$ws.Range('E5').Value = '  +1.67%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '164.55'
$ws.Range('E6').Value = '  -0.92%  '
$ws.Range('D7').Value = '3.812.22'
$ws.Range('E7').Value = '  +0.27%  '
$ws.Range('E8').Value = '  +0.08%  '
$ws.Range('E9').Value = '  -0.14%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.160'
$ws.Range('E10').Value = '  +0.53%  '
$ws.Range('E11').Value = '  -0.48%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.74'
$ws.Range('E12').Value = '  +6.11%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000248'
$ws.Range('E13').Value = '  -0.57%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '35.46'
$ws.Range('E14').Value = '  -1.42%  '
$ws.Range('D15').Value = '4.457.53'
$ws.Range('D16').Value = '3.824.16'
$ws.Range('E16').Value = '  +0.56%  '
$ws.Range('D17').Value = '68.637.11'
$ws.Range('E17').Value = '  +1.23%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '18.11'
$ws.Range('E18').Value = '  -1.40%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.10'
$ws.Range('E19').Value = '  +0.33%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '464.29'
$ws.Range('E21').Value = '  -0.06%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.64'
$ws.Range('E22').Value = '  -1.86%  '
$ws.Range('E23').Value = '  -0.14%  '
$ws.Range('E24').Value = '  +3.59%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '83.79'
$ws.Range('E25').Value = '  +0.62%  '
$ws.Range('E26').Value = '  -0.77%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.11'
$ws.Range('E27').Value = '  -0.20%  '
$ws.Range('E28').Value = '  +0.03%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.98'
$ws.Range('E29').Value = '  -0.30%  '
$ws.Range('D30').Value = '3.963.65'
$ws.Range('E30').Value = '  +0.18%  '
$ws.Range('E31').Value = '  -5.27%  '
$ws.Range('E32').Value = '  +0.07%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '7.25'
$ws.Range('E33').Value = '  -2.22%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '29.04'
$ws.Range('E34').Value = '  -1.28%  '
$ws.Range('E35').Value = '  -0.01%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '9.04'
$ws.Range('E36').Value = '  -0.29%  '
$ws.Range('E37').Value = '  +1.57%  '
$ws.Range('E38').Value = '  +6.30%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.89'
$ws.Range('E39').Value = '  +1.33%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.978'
$ws.Range('E40').Value = '  -1.55%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.15'
$ws.Range('E41').Value = '  -2.22%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.00'
$ws.Range('E42').Value = '  +0.07%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '154.67'
$ws.Range('E44').Value = '  +2.07%  '
$ws.Range('E45').Value = '  -0.48%  '
$ws.Range('B46').Value = 'Arweave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '42.74'
$ws.Range('E46').Value = '  -4.01%  '
$ws.Range('B47').Value = 'ONDO'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.40'
$ws.Range('E47').Value = '  +0.97%  '
$ws.Range('B48').Value = 'OKB'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '46.52'
$ws.Range('E48').Value = '  -2.46%  '
$ws.Range('E49').Value = '  +0.18%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.87'
$ws.Range('E50').Value = '  +1.62%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '378.22'
$ws.Range('E51').Value = '  -2.96%  '
